$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 606, shifting the existing rows 606-630
# down to 607-631 (their content stays unchanged).
$ws.Rows("606:606").Insert()

# Populate the newly inserted row 606 with the new weekly record.
$ws.Range("A606").Value = 10
$ws.Range("B606").Value = "Vega Modelo de Temuco"
$ws.Range("C606").Value = "La Araucanía"
$ws.Range("D606").Value = 44939
$ws.Range("E606").Value = 9
$ws.Range("F606").Value = "Fruta"
$ws.Range("G606").Value = 100108
$ws.Range("H606").Value = "Tropicales y subtropicales"
$ws.Range("I606").Value = 100108005
$ws.Range("J606").Value = "Piña"
$ws.Range("K606").Value = "Caramelo"
$ws.Range("L606").Value = "Segunda"
$ws.Range("M606").Value = 100
$ws.Range("N606").Value = 21000
$ws.Range("O606").Value = 21000
$ws.Range("P606").Value = 21000
$ws.Range("Q606").Value = "$/caja 14 unidades"
$ws.Range("R606").Value = "Ecuador"
$ws.Range("S606").Value = 1500
$ws.Range("T606").Value = 14
